# Generate Report for handback
# Updates the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) timestamps for the
# d189dbdd-... handback row on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-28 09:13:16"
$wsZhCn.Range("G3").Value = "2016-01-28 09:14:06"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-28 09:13:29"
$wsDeDe.Range("G3").Value = "2016-01-28 09:14:26"
